# Update "想去人数" (want-to-go count) values on the sheets that hold the
# 漫展 (comic con) event data: "展览" and "全部类型".
#   F2: 527 -> 528
#   F3: 454 -> 456
#   F4: 14  -> 17

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 528
    $ws.Range("F3").Value = 456
    $ws.Range("F4").Value = 17
}
